$d = $word.ActiveDocument

function Retype-Paragraph($oldText, $newText) {
    # Re-typing the paragraph (via Find & Replace) collapses whatever
    # fragmented runs it used to be split across into a single run that
    # inherits the formatting of the first one - mirroring what Word does
    # whenever a line gets touched while editing.
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

function Fix-SmartQuotes($smartText, $straightText) {
    # Find & Replace auto-smartens straight quotes; put the literal ASCII
    # quotes back via a direct Range.Text assignment (which does not run
    # AutoFormat/AutoCorrect).
    $find = $d.Content.Find
    $find.ClearFormatting()
    $found = $find.Execute($smartText, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($found) {
        $find.Parent.Text = $straightText
    }
}

function Split-RunBoundary($searchText) {
    # Re-typing a paragraph collapses it to one run; toggling Bold on then
    # back off over a sub-range is a formatting no-op that nonetheless
    # forces a fresh run boundary around that sub-range, letting us
    # reproduce the fine-grained run split the target document shows.
    $find = $d.Content.Find
    $find.ClearFormatting()
    $found = $find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($found) {
        $rng = $find.Parent
        $rng.Bold = 1
        $rng.Bold = 0
    }
}

# Paragraphs whose wording is unchanged but which get "retyped" in place so
# the previously run-fragmented text collapses into a single run (mirrors
# what happened when the author touched these lines while editing nearby
# text in Word).
Retype-Paragraph "Análisis de Requerimientos Funcionales y No funcionales" "Análisis de Requerimientos Funcionales y No funcionales"
Retype-Paragraph "El usuario podrá dar de alta una Ficha Técnica de Plantel" "El usuario podrá dar de alta una Ficha Técnica de Plantel"
Retype-Paragraph "El usuario podrá modificar una Ficha Técnica de Plantel." "El usuario podrá modificar una Ficha Técnica de Plantel."
Retype-Paragraph "El Usuario podrá Imprimir una Ficha Técnica de un determinado plantel." "El Usuario podrá Imprimir una Ficha Técnica de un determinado plantel."

Retype-Paragraph "El Sistema de planeación podrá consultar información general y académica de los planteles desde el sistema `"DEO`"." "El Sistema de planeación podrá consultar información general y académica de los planteles desde el sistema `"DEO`"."
Fix-SmartQuotes ([char]0x201C + "DEO" + [char]0x201D) "`"DEO`""

# Actual wording change: the "Ficha Técnica" now *manages* (gestionará) the
# stages, rather than merely storing (almacenará) them. The target keeps
# "gestionará" and "y por tal podrá" each in their own run (mirroring the
# document's existing convention of isolating the verb), so after the bulk
# retype we re-split the run at those two boundaries.
Retype-Paragraph "La ficha Técnica almacenará las etapas de desarrollo del plantel y así agregar, modificar y eliminar dichas etapas." "La ficha Técnica gestionará las etapas de desarrollo del plantel y por tal podrá agregar, modificar y eliminar dichas etapas."
Split-RunBoundary "gestionará"
Split-RunBoundary "y por tal podrá"

Retype-Paragraph "Se podrá subir y descargar evidencias (Fotos) de las etapas de los planteles." "Se podrá subir y descargar evidencias (Fotos) de las etapas de los planteles."

# This paragraph also hosts the "_GoBack" bookmark; retype it first so its
# runs merge, then reposition the bookmark within the merged text below.
Retype-Paragraph "El sistema tendrá acceso restringido a la información por medio " "El sistema tendrá acceso restringido a la información por medio "

Retype-Paragraph "Los usuarios tendrá diferentes roles dentro del sistema." "Los usuarios tendrá diferentes roles dentro del sistema."
Retype-Paragraph "El sistema deberá estar disponible en cualquier momento, para su operación." "El sistema deberá estar disponible en cualquier momento, para su operación."
Retype-Paragraph "El sistema no tendrá restricción en cuanto a los usuario conectados." "El sistema no tendrá restricción en cuanto a los usuario conectados."
Retype-Paragraph "Las impresiones de las Fichas Técnicas deberá adaptarse al tamaño de hoja carta." "Las impresiones de las Fichas Técnicas deberá adaptarse al tamaño de hoja carta."
Retype-Paragraph "El sistema deberá ser programado en JSP de Java." "El sistema deberá ser programado en JSP de Java."
Retype-Paragraph "El Sistema deberá trabajar con una base de datos SQL Server." "El Sistema deberá trabajar con una base de datos SQL Server."

# Move the "_GoBack" bookmark: it used to sit at the very last paragraph
# ("El Sistema deberá trabajar con una base de datos SQL Server.") and now
# marks the most recent edit point, right after "acceso" in
# "El sistema tendrá acceso restringido ...".
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("El sistema tendrá acceso", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackRange = $find.Parent.Duplicate
$goBackRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
